# Generate Report for Handback
# Updates the handoff/handback timestamps for the e4599166-... file (row 3)
# across the per-language sheets, and rolls the newest timestamp up into the
# "Latest HO Xliff Generate Date" column of the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: refresh handoff/handback datetimes for e4599166 row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-19 16:51:57"
$wsZhCn.Range("K3").Value = "2016-08-19 16:52:29"

# --- de-de sheet: refresh handoff/handback datetimes for e4599166 row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-19 16:52:05"
$wsDeDe.Range("K3").Value = "2016-08-19 16:52:36"

# --- Overview sheet: roll up the latest HO Xliff generate date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-19 16:52:05"
